$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'jaah_3'
$ws.Cells.Item(2, 2).Value = 'isophonics_226'
$ws.Cells.Item(2, 3).Value = 0.1527777777777778
$ws.Cells.Item(2, 4).Value = '[[''Ab'', ''Eb'', ''Eb'']]'
$ws.Cells.Item(2, 5).Value = '[[''A'', ''E'', ''E'']]'
$ws.Cells.Item(2, 6).Value = '[(131.54, 134.26)]'
$ws.Cells.Item(2, 7).Value = '[(16.759569, 29.414444)]'
$ws.Cells.Item(2, 8).Value = ""
$ws.Cells.Item(2, 9).Value = ""

$ws.Cells.Item(3, 1).Value = 'schubert-winterreise_36'
$ws.Cells.Item(3, 2).Value = 'isophonics_31'
$ws.Cells.Item(3, 3).Value = 0.1020242914979757
$ws.Cells.Item(3, 4).Value = '[[''G:maj'', ''E:min'', ''A:min/C''], [''G:maj'', ''C:maj/G'', ''G:maj'']]'
$ws.Cells.Item(3, 5).Value = '[[''G'', ''E:min'', ''A:min''], [''G'', ''C'', ''G'']]'
$ws.Cells.Item(3, 6).Value = '[(14.1, 20.24), (21.58, 29.34)]'
$ws.Cells.Item(3, 7).Value = '[(26.115071, 38.793167), (17.581738, 30.271443)]'
$ws.Cells.Item(3, 8).Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Cells.Item(3, 9).Value = ""

$ws.Cells.Item(4, 1).Value = 'isophonics_235'
$ws.Cells.Item(4, 2).Value = 'schubert-winterreise_72'
$ws.Cells.Item(4, 3).Value = 0.1458333333333333
$ws.Cells.Item(4, 4).Value = '[[''Bb'', ''F'', ''Bb''], [''F'', ''Bb'', ''Eb:7'']]'
$ws.Cells.Item(4, 5).Value = '[[''C:maj'', ''G:maj'', ''C:maj''], [''G:maj'', ''C:maj'', ''F:7'']]'
$ws.Cells.Item(4, 6).Value = '[(26.041712, 37.10602), (27.910918, 40.797993)]'
$ws.Cells.Item(4, 7).Value = '[(8.34, 10.04), (12.24, 15.4)]'
$ws.Cells.Item(4, 8).Value = ""
$ws.Cells.Item(4, 9).Value = ""

$ws.Cells.Item(5, 1).Value = 'isophonics_296'
$ws.Cells.Item(5, 2).Value = 'schubert-winterreise_9'
$ws.Cells.Item(5, 3).Value = 0.1001011122345804
$ws.Cells.Item(5, 4).Value = '[[''D:maj'', ''A/3'', ''D/7'']]'
$ws.Cells.Item(5, 5).Value = '[[''F:maj'', ''C:maj'', ''F:maj'']]'
$ws.Cells.Item(5, 6).Value = '[(59.78, 64.294)]'
$ws.Cells.Item(5, 7).Value = '[(46.6, 49.04)]'
$ws.Cells.Item(5, 8).Value = ""
$ws.Cells.Item(5, 9).Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

$ws.Cells.Item(6, 1).Value = 'schubert-winterreise_40'
$ws.Cells.Item(6, 2).Value = 'schubert-winterreise_147'
$ws.Cells.Item(6, 3).Value = 0.5397727272727273
$ws.Cells.Item(6, 4).Value = '[[''D:maj'', ''A:7'', ''D:maj'', ''A:7'', ''D:maj'']]'
$ws.Cells.Item(6, 5).Value = '[[''A:maj/E'', ''E:7'', ''A:maj'', ''E:7'', ''A:maj'']]'
$ws.Cells.Item(6, 6).Value = '[(27.0, 49.86)]'
$ws.Cells.Item(6, 7).Value = '[(19.78, 25.82)]'
$ws.Cells.Item(6, 8).Value = ""
$ws.Cells.Item(6, 9).Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'

$ws.Cells.Item(7, 1).Value = 'isophonics_227'
$ws.Cells.Item(7, 2).Value = 'schubert-winterreise_163'
$ws.Cells.Item(7, 3).Value = 0.13125
$ws.Cells.Item(7, 4).Value = '[[''E:7'', ''A:min'', ''A:min/b3'']]'
$ws.Cells.Item(7, 5).Value = '[[''C:7'', ''F:min/C'', ''F:min'']]'
$ws.Cells.Item(7, 6).Value = '[(5.67204, 8.841564)]'
$ws.Cells.Item(7, 7).Value = '[(32.32, 33.88)]'
$ws.Cells.Item(7, 8).Value = ""
$ws.Cells.Item(7, 9).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

$ws.Cells.Item(8, 1).Value = 'isophonics_57'
$ws.Cells.Item(8, 2).Value = 'isophonics_1'
$ws.Cells.Item(8, 3).Value = 0.1039136302294197
$ws.Cells.Item(8, 4).Value = '[[''G'', ''C'', ''G'', ''C'']]'
$ws.Cells.Item(8, 5).Value = '[[''Eb'', ''Ab/5'', ''Eb'', ''Ab/5'']]'
$ws.Cells.Item(8, 6).Value = '[(135.68755, 139.228594)]'
$ws.Cells.Item(8, 7).Value = '[(17.016, 24.727)]'
$ws.Cells.Item(8, 8).Value = ""
$ws.Cells.Item(8, 9).Value = ""

$ws.Cells.Item(9, 1).Value = 'schubert-winterreise_186'
$ws.Cells.Item(9, 2).Value = 'jaah_87'
$ws.Cells.Item(9, 3).Value = 0.1525641025641026
$ws.Cells.Item(9, 4).Value = '[[''F:min'', ''C:maj'', ''F:min'', ''G:maj/B'']]'
$ws.Cells.Item(9, 5).Value = '[[''F:min'', ''C'', ''F:min'', ''G'']]'
$ws.Cells.Item(9, 6).Value = '[(15.04, 20.04)]'
$ws.Cells.Item(9, 7).Value = '[(3.4, 9.08)]'
$ws.Cells.Item(9, 8).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Cells.Item(9, 9).Value = ""

$ws.Cells.Item(10, 1).Value = 'isophonics_282'
$ws.Cells.Item(10, 2).Value = 'isophonics_200'
$ws.Cells.Item(10, 3).Value = 0.1652173913043478
$ws.Cells.Item(10, 4).Value = '[[''Eb'', ''Bb'', ''F''], [''Eb'', ''F'', ''Bb'']]'
$ws.Cells.Item(10, 5).Value = '[[''D'', ''A'', ''E''], [''D'', ''E'', ''A'']]'
$ws.Cells.Item(10, 6).Value = '[(29.719863, 35.722222), (44.917324, 48.272607)]'
$ws.Cells.Item(10, 7).Value = '[(8.085475, 18.139715), (60.040123, 71.429511)]'
$ws.Cells.Item(10, 8).Value = 'spotify:track:3Am0IbOxmvlSXro7N5iSfZ'
$ws.Cells.Item(10, 9).Value = ""

$ws.Cells.Item(11, 1).Value = 'schubert-winterreise_113'
$ws.Cells.Item(11, 2).Value = 'schubert-winterreise_161'
$ws.Cells.Item(11, 3).Value = 0.2708333333333333
$ws.Cells.Item(11, 4).Value = '[[''F:maj'', ''C:7'', ''F:maj'', ''C:7'', ''F:maj'']]'
$ws.Cells.Item(11, 5).Value = '[[''C:maj'', ''G:7'', ''C:maj'', ''G:7'', ''C:maj'']]'
$ws.Cells.Item(11, 6).Value = '[(59.9, 69.84)]'
$ws.Cells.Item(11, 7).Value = '[(1.6, 12.54)]'
$ws.Cells.Item(11, 8).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Cells.Item(11, 9).Value = ""

$ws.Cells.Item(12, 1).Value = 'schubert-winterreise_14'
$ws.Cells.Item(12, 2).Value = 'isophonics_291'
$ws.Cells.Item(12, 3).Value = 0.1350574712643678
$ws.Cells.Item(12, 4).Value = '[[''D:maj'', ''G:maj'', ''D:maj/F#'']]'
$ws.Cells.Item(12, 5).Value = '[[''D'', ''G'', ''D'']]'
$ws.Cells.Item(12, 6).Value = '[(83.34, 87.64)]'
$ws.Cells.Item(12, 7).Value = '[(5.20815, 11.03585)]'
$ws.Cells.Item(12, 8).Value = ""
$ws.Cells.Item(12, 9).Value = 'spotify:track:06ypiqmILMdVeaiErMFA91'

$ws.Cells.Item(13, 1).Value = 'schubert-winterreise_152'
$ws.Cells.Item(13, 2).Value = 'schubert-winterreise_15'
$ws.Cells.Item(13, 3).Value = 0.1348837209302326
$ws.Cells.Item(13, 4).Value = '[[''E:7'', ''A:min'', ''A:min'']]'
$ws.Cells.Item(13, 5).Value = '[[''E:7'', ''A:min'', ''A:min'']]'
$ws.Cells.Item(13, 6).Value = '[(13.66, 17.9)]'
$ws.Cells.Item(13, 7).Value = '[(7.36, 16.52)]'
$ws.Cells.Item(13, 8).Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'
$ws.Cells.Item(13, 9).Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'

$ws.Cells.Item(14, 1).Value = 'isophonics_248'
$ws.Cells.Item(14, 2).Value = 'isophonics_283'
$ws.Cells.Item(14, 3).Value = 0.15
$ws.Cells.Item(14, 4).Value = '[[''B'', ''E/5'', ''F#/4'', ''B'']]'
$ws.Cells.Item(14, 5).Value = '[[''G'', ''C'', ''D'', ''G'']]'
$ws.Cells.Item(14, 6).Value = '[(123.652, 131.396)]'
$ws.Cells.Item(14, 7).Value = '[(63.629047, 74.101247)]'
$ws.Cells.Item(14, 8).Value = ""
$ws.Cells.Item(14, 9).Value = 'spotify:track:3tGhRLgcCP6SIZU3tbGl7l'

$ws.Cells.Item(15, 1).Value = 'isophonics_51'
$ws.Cells.Item(15, 2).Value = 'isophonics_235'
$ws.Cells.Item(15, 3).Value = 0.1458333333333333
$ws.Cells.Item(15, 4).Value = '[[''F'', ''G'', ''C'']]'
$ws.Cells.Item(15, 5).Value = '[[''Eb'', ''F'', ''Bb'']]'
$ws.Cells.Item(15, 6).Value = '[(25.687641, 32.084739)]'
$ws.Cells.Item(15, 7).Value = '[(59.292687, 74.118628)]'
$ws.Cells.Item(15, 8).Value = ""
$ws.Cells.Item(15, 9).Value = ""

$ws.Cells.Item(16, 1).Value = 'schubert-winterreise_155'
$ws.Cells.Item(16, 2).Value = 'schubert-winterreise_104'
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = '[[''C:maj'', ''G:7'', ''C:maj'', ''G:7'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''F:7'', ''A#:maj'', ''D:7/A'', ''G#:(3,5)'', ''C:maj'', ''A:(3,5,b7,b9)'', ''D:min/G'', ''F:maj/G'', ''G:7'', ''C:maj'', ''C:min'']]'
$ws.Cells.Item(16, 5).Value = '[[''C:maj'', ''G:7'', ''C:maj'', ''G:7'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''F:7'', ''A#:maj'', ''D:7/A'', ''G#:(3,5)'', ''C:maj'', ''A:(3,5,b7,b9)'', ''D:min/G'', ''F:maj/G'', ''G:7'', ''C:maj'', ''C:min'']]'
$ws.Cells.Item(16, 6).Value = '[(1.46, 39.82)]'
$ws.Cells.Item(16, 7).Value = '[(0.24, 35.1)]'
$ws.Cells.Item(16, 8).Value = ""
$ws.Cells.Item(16, 9).Value = ""

$ws.Cells.Item(17, 1).Value = 'isophonics_204'
$ws.Cells.Item(17, 2).Value = 'isophonics_128'
$ws.Cells.Item(17, 3).Value = 0.1366396761133603
$ws.Cells.Item(17, 4).Value = '[[''D'', ''A'', ''D''], [''A'', ''D'', ''A'']]'
$ws.Cells.Item(17, 5).Value = '[[''F'', ''C'', ''F''], [''C'', ''F'', ''C'']]'
$ws.Cells.Item(17, 6).Value = '[(32.874217, 37.657528), (33.698526, 40.885102)]'
$ws.Cells.Item(17, 7).Value = '[(12.376258, 19.597664), (10.634761, 15.986961)]'
$ws.Cells.Item(17, 8).Value = ""
$ws.Cells.Item(17, 9).Value = ""
